$d = $word.ActiveDocument

# Locate the "Edison Achalma" byline paragraph (uses the "Author" style,
# directly under the document title) and insert a new "Author"-styled
# paragraph right after it containing the author's affiliation.
$targetIndex = 0
$count = 0
foreach ($p in $d.Paragraphs) {
    $count = $count + 1
    $styleName = $p.Style.NameLocal
    $text = $p.Range.Text.Trim()
    if ($styleName -eq "Author" -and $text -eq "Edison Achalma") {
        $targetIndex = $count
        break
    }
}

if ($targetIndex -gt 0) {
    # The paragraph immediately following the byline (a blank "Body Text"
    # paragraph) is where the new paragraph needs to land, so split right
    # before its content -- this leaves the "Edison Achalma" paragraph
    # completely untouched and creates a fresh paragraph in between.
    $nextPara = $d.Paragraphs($targetIndex + 1)
    $nextPara.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs($targetIndex + 1)
    $newPara.Style = $d.Styles("Author")
    $newPara.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
}
